# correção nos dados e inicio da analise PNAD 2009
#
# The sheet had two "section header" rows that carried only a label in
# column A and no data in B:F:
#   row 5 -> "situação do domicílio"
#   row 8 -> "grandes regiões e unidades da federação"
# Both are removed (EntireRow delete), which shifts every row below up
# (rows 6-7 move to 5-6, rows 9-39 move to 7-37). The now-unused shared
# strings for those two headers (and the "unnamed: 1_level_1" /
# "unnamed: 5_level_1" placeholders) are dropped automatically once
# nothing references them anymore.
#
# The row-2 sub-header cells that used to read "unnamed: 1_level_1" (B2)
# and "unnamed: 5_level_1" (F2) are corrected to read "total", matching
# the existing "total" label already used in C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "grandes regiões e unidades da federação" header row first
# (row 8) then the "situação do domicílio" header row (row 5) - deleting
# the lower row first keeps the row-5 index valid for the second delete.
$ws.Rows(8).Delete()
$ws.Rows(5).Delete()

# Fix the row-2 sub-headers.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
